$wb = $excel.ActiveWorkbook

# "User" sheet: field in row 5 renamed from "username" to "name"
$wsUser = $wb.Worksheets.Item("User")
$wsUser.Range("A5").Value = "name"

# "SharedImage" sheet: new row documenting the image-data-format reference link
$wsShared = $wb.Worksheets.Item("SharedImage")
$wsShared.Range("A7").Value = "https://tpiros.dev/blog/image-management-via-graphql/"

# Reproduce per-sheet selection state
$wsUser.Activate()
$wsUser.Range("A8").Select()

$wsShared.Activate()
$wsShared.Range("C17").Select()

$wsCreated = $wb.Worksheets.Item("CreatedImage")
$wsCreated.Activate()
$wsCreated.Range("C37").Select()
